$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "B2" = 0.0006408296065709695
    "C2" = 0.04071648406533734
    "D2" = 0.1494219747398047
    "E2" = 0.4942365360607697
    "G2" = 0.6850158244724827

    "B3" = 3.286832544864788
    "C3" = 1.655778082260271
    "D3" = 3.537761648806719
    "E3" = 0.4942365360607697
    "G3" = 8.974608811992548

    "B4" = 3.286832544864788
    "C4" = 1.655778082260271
    "D4" = 0.1494219747398047
    "E4" = 0.4942365360607697
    "G4" = 5.586269137925634

    "B5" = 0.04271373187048222
    "C5" = 0.04071648406533734
    "D5" = 0.1494219747398047
    "E5" = 0.4942365360607697
    "G5" = 0.7270887267363939

    "B6" = 3.286832544864788
    "C6" = 1.655778082260271
    "D6" = 0.7527432677738641
    "E6" = 0.4942365360607697
    "G6" = 6.189590430959694
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
